$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "26.309.62"
Set-TextValue $ws.Range("E2") "  +0.48%  "
Set-TextValue $ws.Range("D3") "1.595.71"
Set-TextValue $ws.Range("E3") "  +0.28%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "211.73"
Set-TextValue $ws.Range("E5") "  -0.12%  "
Set-TextValue $ws.Range("D9") "0.0605"
Set-TextValue $ws.Range("E9") "  +0.01%  "
Set-TextValue $ws.Range("E10") "  +0.33%  "
Set-TextValue $ws.Range("D11") "0.0856"
Set-TextValue $ws.Range("E11") "  +1.25%  "
Set-TextValue $ws.Range("D12") "1.819.95"
Set-TextValue $ws.Range("E12") "  +0.27%  "
Set-TextValue $ws.Range("D13") "1.608.37"
Set-TextValue $ws.Range("E13") "  +1.00%  "
Set-TextValue $ws.Range("E14") "  -0.62%  "
Set-TextValue $ws.Range("E15") "  -1.07%  "
Set-TextValue $ws.Range("D16") "63.44"
Set-TextValue $ws.Range("E16") "  -0.24%  "
Set-TextValue $ws.Range("D17") "26.307.44"
Set-TextValue $ws.Range("E17") "  +0.44%  "
Set-TextValue $ws.Range("D18") "229.74"
Set-TextValue $ws.Range("E18") "  +7.55%  "
Set-TextValue $ws.Range("B19") "Chainlink"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D19") "7.64"
Set-TextValue $ws.Range("E19") "  +4.04%  "
Set-TextValue $ws.Range("B20") "ShibaInu"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D20") "0.0₃0721"
Set-TextValue $ws.Range("E20") "  -0.55%  "
Set-TextValue $ws.Range("E21") "  -0.05%  "
Set-TextValue $ws.Range("E22") "  -0.24%  "
Set-TextValue $ws.Range("E23") "  +2.77%  "
Set-TextValue $ws.Range("E24") "  -1.20%  "
Set-TextValue $ws.Range("D25") "146.47"
Set-TextValue $ws.Range("E25") "  +1.15%  "
Set-TextValue $ws.Range("E26") "  +0.08%  "
Set-TextValue $ws.Range("E27") "  +0.05%  "
Set-TextValue $ws.Range("E28") "  +0.31%  "
Set-TextValue $ws.Range("D29") "15.36"
Set-TextValue $ws.Range("E29") "  +1.85%  "
Set-TextValue $ws.Range("E30") "  +0.11%  "
Set-TextValue $ws.Range("E31") "  -0.27%  "
Set-TextValue $ws.Range("D32") "1.497.27"
Set-TextValue $ws.Range("E32") "  +5.32%  "
Set-TextValue $ws.Range("E33") "  +1.32%  "
Set-TextValue $ws.Range("D34") "2.93"
Set-TextValue $ws.Range("E34") "  -0.95%  "
Set-TextValue $ws.Range("E35") "  -0.25%  "
Set-TextValue $ws.Range("E36") "  +0.61%  "
Set-TextValue $ws.Range("D37") "0.568"
Set-TextValue $ws.Range("E37") "  -3.14%  "
Set-TextValue $ws.Range("E38") "  -0.72%  "
Set-TextValue $ws.Range("D39") "0.817"
Set-TextValue $ws.Range("E39") "  -0.66%  "
Set-TextValue $ws.Range("E40") "  -2.09%  "
Set-TextValue $ws.Range("E41") "  +0.08%  "
Set-TextValue $ws.Range("E42") "  +1.87%  "
Set-TextValue $ws.Range("E43") "  -3.80%  "
Set-TextValue $ws.Range("D44") "1.733.06"
Set-TextValue $ws.Range("E44") "  +0.38%  "
Set-TextValue $ws.Range("E45") "  -1.12%  "
Set-TextValue $ws.Range("D46") "60.61"
Set-TextValue $ws.Range("E46") "  -0.54%  "
Set-TextValue $ws.Range("E47") "  +1.73%  "
Set-TextValue $ws.Range("E49") "  -0.26%  "
Set-TextValue $ws.Range("E50") "  -0.10%  "
Set-TextValue $ws.Range("D51") "0.999"
Set-TextValue $ws.Range("E51") "  +0.10%  "
